# Refresh the crypto price/volume table to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subscript-3 glyph (U+2083) used in the PEPE price cell; build it via the
# format operator rather than string concatenation, since "0.0" + [char] is
# evaluated as numeric addition, not concatenation, in this host.
$sub3 = [char]0x2083

$ws.Range('D2').Value = '63.414.52'
$ws.Range('E2').Value = '  +2.04%  '
$ws.Range('D3').Value = '3.172.35'
$ws.Range('E3').Value = '  +0.20%  '
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').Value = "`'602.11"
$ws.Range('E5').Value = '  +2.35%  '
$ws.Range('D6').Value = "`'136.08"
$ws.Range('E6').Value = '  +1.27%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '3.169.35'
$ws.Range('D9').Value = "`'0.513"
$ws.Range('E9').Value = '  +2.52%  '
$ws.Range('E10').Value = '  +1.14%  '
$ws.Range('D11').Value = "`'5.37"
$ws.Range('E11').Value = '  +2.64%  '
$ws.Range('D12').Value = "`'0.455"
$ws.Range('E12').Value = '  +1.62%  '
$ws.Range('E13').Value = '  +2.81%  '
$ws.Range('D14').Value = "`'34.86"
$ws.Range('E14').Value = '  +5.50%  '
$ws.Range('D15').Value = '3.695.57'
$ws.Range('E15').Value = '  -0.15%  '
$ws.Range('E16').Value = '  +0.52%  '
$ws.Range('D17').Value = '3.171.42'
$ws.Range('E17').Value = '  -0.19%  '
$ws.Range('D18').Value = '63.395.63'
$ws.Range('E18').Value = '  +1.92%  '
$ws.Range('D19').Value = "`'6.59"
$ws.Range('E19').Value = '  +0.91%  '
$ws.Range('D20').Value = "`'462.47"
$ws.Range('E20').Value = '  +1.44%  '
$ws.Range('E21').Value = '  +1.06%  '
$ws.Range('E22').Value = '  -0.19%  '
$ws.Range('E23').Value = '  +1.52%  '
$ws.Range('D24').Value = "`'13.30"
$ws.Range('E24').Value = '  +0.65%  '
$ws.Range('D25').Value = "`'83.20"
$ws.Range('E25').Value = '  +1.31%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('E27').Value = '  +1.17%  '
$ws.Range('E28').Value = '  -0.20%  '
$ws.Range('E29').Value = '  +2.91%  '
$ws.Range('E30').Value = '  -0.99%  '
$ws.Range('D31').Value = "`'7.72"
$ws.Range('E31').Value = '  -0.77%  '
$ws.Range('D32').Value = "`'27.16"
$ws.Range('E32').Value = '  +0.26%  '
$ws.Range('E33').Value = '  -0.70%  '
$ws.Range('D34').Value = "`'2.44"
$ws.Range('E34').Value = '  +2.49%  '
$ws.Range('D35').Value = "`'1.02"
$ws.Range('E35').Value = '  -1.08%  '
$ws.Range('E36').Value = '  +2.57%  '
$ws.Range('B37').Value = 'PEPE'
$ws.Range('C37').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D37').Value = ("0.0{0}0733" -f $sub3)
$ws.Range('E37').Value = '  +7.63%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').Value = "`'51.23"
$ws.Range('E38').Value = '  +0.44%  '
$ws.Range('D39').Value = "`'0.0390"
$ws.Range('E39').Value = '  +1.45%  '
$ws.Range('E40').Value = '  +2.16%  '
$ws.Range('E41').Value = '  +0.96%  '
$ws.Range('D42').Value = "`'2.64"
$ws.Range('E42').Value = '  +0.26%  '
$ws.Range('D43').Value = "`'394.18"
$ws.Range('E43').Value = '  -3.46%  '
$ws.Range('D44').Value = '2.808.83'
$ws.Range('E44').Value = '  -4.56%  '
$ws.Range('D45').Value = "`'0.251"
$ws.Range('E45').Value = '  +1.24%  '
$ws.Range('D46').Value = "`'36.07"
$ws.Range('E46').Value = '  +1.21%  '
$ws.Range('E48').Value = '  -0.20%  '
$ws.Range('D49').Value = "`'126.27"
$ws.Range('E49').Value = '  +2.27%  '
$ws.Range('D50').Value = "`'25.19"
$ws.Range('E50').Value = '  -0.61%  '
$ws.Range('E51').Value = '  +0.75%  '
